$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format before writing, so numeric-looking
# strings like "378.91" are stored as literal text instead of being coerced
# to a floating point number by Excel's type inference on Range.Value.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.787.64"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "2.975.32"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "378.91"
$ws.Range("E5").Value = "  +7.30%  "
$ws.Range("D6").Value = "104.96"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").Value = "37.51"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "18.68"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").Value = "3.442.25"
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "2.985.32"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").Value = "0.960"
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").Value = "51.853.65"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "3.47"
$ws.Range("E19").Value = "  +4.23%  "
$ws.Range("D20").Value = "7.43"
$ws.Range("E20").Value = "  +1.79%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "0.0₃0960"
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("D23").Value = "68.81"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "263.74"
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("E25").Value = "  +4.55%  "
$ws.Range("D26").Value = "7.46"
$ws.Range("E26").Value = "  +19.17%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.170"
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "4.16"
$ws.Range("E28").Value = "  -3.99%  "
$ws.Range("E29").Value = "  +3.84%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "26.07"
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("E32").Value = "  -2.93%  "
$ws.Range("D33").Value = "9.94"
$ws.Range("D34").Value = "51.64"
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("D35").Value = "34.60"
$ws.Range("E35").Value = "  -2.37%  "
$ws.Range("E36").Value = "  -4.15%  "
$ws.Range("D37").Value = "0.0437"
$ws.Range("E37").Value = "  +2.46%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").Value = "3.08"
$ws.Range("E39").Value = "  -4.54%  "
$ws.Range("D40").Value = "17.46"
$ws.Range("E40").Value = "  +1.18%  "
$ws.Range("E41").Value = "  -5.81%  "
$ws.Range("D42").Value = "1.86"
$ws.Range("E42").Value = "  -2.43%  "
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("D44").Value = "124.02"
$ws.Range("E44").Value = "  +2.67%  "
$ws.Range("D45").Value = "22.26"
$ws.Range("E45").Value = "  -2.75%  "
$ws.Range("E46").Value = "  +18.68%  "
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("D48").Value = "2.038.01"
$ws.Range("E48").Value = "  -2.89%  "
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("D50").Value = "3.24"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").Value = "0.0335"
$ws.Range("E51").Value = "  +4.99%  "

# Restore the default (unstyled) cell style now that the text values are
# committed, so the cells end up with no explicit style index, matching
# the original workbook's formatting.
$ws.Range("D2:D51").Style = "Normal"
